$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 144.6
$ws.Range("I9").Value = 180
$ws.Range("J9").Value = 113.625
$ws.Range("K9").Value = 180
$ws.Range("L9").Value = 113.625
$ws.Range("M9").Value = -11
$ws.Range("N9").Value = -451.625
$ws.Range("H17").Value = 288.3
$ws.Range("J17").Value = 288.3
$ws.Range("L17").Value = 864.9000000000001
$ws.Range("N17").Value = -1200.9
$ws.Range("H42").Value = 596.94116
$ws.Range("I42").Value = 123.333336
$ws.Range("J42").Value = 1733.6
$ws.Range("K42").Value = 370.000008
$ws.Range("L42").Value = 5200.799999999999
$ws.Range("M42").Value = -140.000008
$ws.Range("N42").Value = -5660.799999999999
$ws.Range("H61").Value = 41666764
$ws.Range("I61").Value = 102.5
$ws.Range("J61").Value = 166666750
$ws.Range("K61").Value = 307.5
$ws.Range("L61").Value = 500000250
$ws.Range("M61").Value = -135.5
$ws.Range("N61").Value = -500000594
$ws.Range("H80").Value = 512.3214
$ws.Range("I80").Value = 551.5833
$ws.Range("J80").Value = 482.875
$ws.Range("K80").Value = 1654.7499
$ws.Range("L80").Value = 1448.625
$ws.Range("M80").Value = -656.7499
$ws.Range("N80").Value = -3444.625
$ws.Range("H83").Value = 512.3214
$ws.Range("I83").Value = 551.5833
$ws.Range("J83").Value = 482.875
$ws.Range("K83").Value = 4964.2497
$ws.Range("L83").Value = 4345.875
$ws.Range("M83").Value = 27.7502999999997
$ws.Range("N83").Value = -14329.875
$ws.Range("H86").Value = 1754.7
$ws.Range("I86").Value = 1747.1177
$ws.Range("K86").Value = 1747.1177
$ws.Range("M86").Value = -624.1177
$ws.Range("H89").Value = 1754.7
$ws.Range("I89").Value = 1747.1177
$ws.Range("K89").Value = 8735.5885
$ws.Range("M89").Value = -3119.5885
$ws.Range("H92").Value = 1241.5834
$ws.Range("I92").Value = 285.57144
$ws.Range("J92").Value = 2580
$ws.Range("K92").Value = 285.57144
$ws.Range("L92").Value = 2580
$ws.Range("M92").Value = 962.4285600000001
$ws.Range("N92").Value = -5076
$ws.Range("H103").Value = 1155
$ws.Range("I103").Value = 721.8
$ws.Range("K103").Value = 2165.4
$ws.Range("M103").Value = -1579.4
$ws.Range("H106").Value = 2835
$ws.Range("I106").Value = 2752.5
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2752.5
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -2121.5
$ws.Range("N106").Value = -4262
$ws.Range("H107").Value = 479.2
$ws.Range("I107").Value = 460.82352
$ws.Range("J107").Value = 583.3333
$ws.Range("K107").Value = 460.82352
$ws.Range("L107").Value = 583.3333
$ws.Range("M107").Value = 1459.17648
$ws.Range("N107").Value = -4423.3333
$ws.Range("H112").Value = 2439.1
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 2654.5557
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 7963.6671
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -10179.6671
$ws.Range("H116").Value = 3946.9333
$ws.Range("I116").Value = 2872.8572
$ws.Range("J116").Value = 4273.826
$ws.Range("K116").Value = 2872.8572
$ws.Range("L116").Value = 4273.826
$ws.Range("M116").Value = 569.1428000000001
$ws.Range("N116").Value = -11157.826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21742168
$ws.Range("I32").Value = 23257322
$ws.Range("K32").Value = 23257322
$ws.Range("M32").Value = -23257035
$ws.Range("H45").Value = 2311.7
$ws.Range("I45").Value = 1868.8462
$ws.Range("K45").Value = 1868.8462
$ws.Range("M45").Value = -1491.8462
$ws.Range("H74").Value = 4732.5356
$ws.Range("I74").Value = 6105.1055
$ws.Range("J74").Value = 1834.8889
$ws.Range("K74").Value = 6105.1055
$ws.Range("L74").Value = 1834.8889
$ws.Range("M74").Value = -5231.1055
$ws.Range("N74").Value = -3582.8889
$ws.Range("H77").Value = 4732.5356
$ws.Range("I77").Value = 6105.1055
$ws.Range("J77").Value = 1834.8889
$ws.Range("K77").Value = 30525.5275
$ws.Range("L77").Value = 9174.4445
$ws.Range("M77").Value = -26157.5275
$ws.Range("N77").Value = -17910.4445
$ws.Range("H102").Value = 3397.5
$ws.Range("I102").Value = 2905
$ws.Range("K102").Value = 2905
$ws.Range("M102").Value = -1283
$ws.Range("H110").Value = 2953.0454
$ws.Range("I110").Value = 2314.5789
$ws.Range("J110").Value = 6996.6665
$ws.Range("K110").Value = 2314.5789
$ws.Range("L110").Value = 6996.6665
$ws.Range("M110").Value = -269.5789
$ws.Range("N110").Value = -11086.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 862.5263
$ws.Range("I94").Value = 797.2
$ws.Range("J94").Value = 1107.5
$ws.Range("K94").Value = 797.2
$ws.Range("L94").Value = 1107.5
$ws.Range("M94").Value = -346.2
$ws.Range("N94").Value = -2009.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2933.3333
$ws.Range("I99").Value = 2920
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2920
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1422
$ws.Range("N99").Value = -5996
$ws.Range("H122").Value = 1307.6451
$ws.Range("I122").Value = 1181.421
$ws.Range("J122").Value = 1507.5
$ws.Range("K122").Value = 3544.263
$ws.Range("L122").Value = 4522.5
$ws.Range("M122").Value = -1094.263
$ws.Range("N122").Value = -9422.5
$ws.Range("H126").Value = 2933.3333
$ws.Range("I126").Value = 2920
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8760
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6290
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 246.54546
$ws.Range("I6").Value = 71.2
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 213.6
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = -100.6
$ws.Range("N6").Value = -6226
$ws.Range("H50").Value = 828.8
$ws.Range("I50").Value = 48
$ws.Range("J50").Value = 2000
$ws.Range("K50").Value = 144
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = 337
$ws.Range("N50").Value = -6962
$ws.Range("H52").Value = 416.25
$ws.Range("J52").Value = 416.25
$ws.Range("L52").Value = 1248.75
$ws.Range("N52").Value = -1780.75
$ws.Range("H53").Value = 828.8
$ws.Range("I53").Value = 48
$ws.Range("J53").Value = 2000
$ws.Range("K53").Value = 144
$ws.Range("L53").Value = 6000
$ws.Range("M53").Value = 337
$ws.Range("N53").Value = -6962
$ws.Range("H93").Value = 5000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 15000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -18744
$ws.Range("H122").Value = 744.6316
$ws.Range("I122").Value = 631.7273
$ws.Range("J122").Value = 899.875
$ws.Range("K122").Value = 5685.545700000001
$ws.Range("L122").Value = 8098.875
$ws.Range("M122").Value = -3235.545700000001
$ws.Range("N122").Value = -12998.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 4500
$ws.Range("J47").Value = 4500
$ws.Range("L47").Value = 4500
$ws.Range("N47").Value = -5636
$ws.Range("H70").Value = 5516.846
$ws.Range("I70").Value = 5850.241
$ws.Range("J70").Value = 4550
$ws.Range("K70").Value = 5850.241
$ws.Range("L70").Value = 4550
$ws.Range("M70").Value = -5580.241
$ws.Range("N70").Value = -5090
$ws.Range("H73").Value = 5516.846
$ws.Range("I73").Value = 5850.241
$ws.Range("J73").Value = 4550
$ws.Range("K73").Value = 5850.241
$ws.Range("L73").Value = 4550
$ws.Range("M73").Value = -4914.241
$ws.Range("N73").Value = -6422
$ws.Range("H139").Value = 32000
$ws.Range("J139").Value = 32000
$ws.Range("L139").Value = 32000
$ws.Range("N139").Value = -42280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1429686
$ws.Range("I22").Value = 2500575
$ws.Range("J22").Value = 1834
$ws.Range("K22").Value = 2500575
$ws.Range("L22").Value = 1834
$ws.Range("M22").Value = -2500280
$ws.Range("N22").Value = -2424
$ws.Range("H27").Value = 1429686
$ws.Range("I27").Value = 2500575
$ws.Range("J27").Value = 1834
$ws.Range("K27").Value = 2500575
$ws.Range("L27").Value = 1834
$ws.Range("M27").Value = -2500468
$ws.Range("N27").Value = -2048
$ws.Range("H46").Value = 1863
